$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3500.5
$ws.Range("J62").Value = 5006
$ws.Range("L62").Value = 5006
$ws.Range("N62").Value = -6254
$ws.Range("H65").Value = 3500.5
$ws.Range("J65").Value = 5006
$ws.Range("L65").Value = 25030
$ws.Range("N65").Value = -31270
$ws.Range("H86").Value = 1738
$ws.Range("I86").Value = 1497.5
$ws.Range("J86").Value = 1858.25
$ws.Range("K86").Value = 1497.5
$ws.Range("L86").Value = 1858.25
$ws.Range("M86").Value = -374.5
$ws.Range("N86").Value = -4104.25
$ws.Range("H89").Value = 1738
$ws.Range("I89").Value = 1497.5
$ws.Range("J89").Value = 1858.25
$ws.Range("K89").Value = 7487.5
$ws.Range("L89").Value = 9291.25
$ws.Range("M89").Value = -1871.5
$ws.Range("N89").Value = -20523.25
$ws.Range("H132").Value = 1430171.4
$ws.Range("I132").Value = 1638.8
$ws.Range("J132").Value = 5001503
$ws.Range("K132").Value = 4916.4
$ws.Range("L132").Value = 15004509
$ws.Range("M132").Value = -2386.4
$ws.Range("N132").Value = -15009569

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6028.931
$ws.Range("I32").Value = 3139.818
$ws.Range("J32").Value = 58996
$ws.Range("K32").Value = 3139.818
$ws.Range("L32").Value = 58996
$ws.Range("M32").Value = -2852.818
$ws.Range("N32").Value = -59570
$ws.Range("H61").Value = 2425.5151
$ws.Range("I61").Value = 1790.96
$ws.Range("K61").Value = 1790.96
$ws.Range("M61").Value = -1578.96
$ws.Range("H74").Value = 875.375
$ws.Range("I74").Value = 821.8570999999999
$ws.Range("J74").Value = 1250
$ws.Range("K74").Value = 821.8570999999999
$ws.Range("L74").Value = 1250
$ws.Range("M74").Value = 52.14290000000005
$ws.Range("N74").Value = -2998
$ws.Range("H77").Value = 875.375
$ws.Range("I77").Value = 821.8570999999999
$ws.Range("J77").Value = 1250
$ws.Range("K77").Value = 4109.2855
$ws.Range("L77").Value = 6250
$ws.Range("M77").Value = 258.7145
$ws.Range("N77").Value = -14986
$ws.Range("H132").Value = 1899.5938
$ws.Range("I132").Value = 1947.4138
$ws.Range("J132").Value = 1437.3334
$ws.Range("K132").Value = 5842.2414
$ws.Range("L132").Value = 4312.0002
$ws.Range("M132").Value = -3312.2414
$ws.Range("N132").Value = -9372.0002
$ws.Range("H134").Value = 88750
$ws.Range("J134").Value = 88750
$ws.Range("L134").Value = 88750
$ws.Range("N134").Value = -98890
$ws.Range("H136").Value = 2425.5151
$ws.Range("I136").Value = 1790.96
$ws.Range("K136").Value = 5372.88
$ws.Range("M136").Value = -2822.88
$ws.Range("H141").Value = 58333
$ws.Range("J141").Value = 67499.5
$ws.Range("L141").Value = 67499.5
$ws.Range("N141").Value = -77859.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 835.8333
$ws.Range("I20").Value = 910.0833
$ws.Range("J20").Value = 687.3333
$ws.Range("K20").Value = 910.0833
$ws.Range("L20").Value = 687.3333
$ws.Range("M20").Value = -663.0833
$ws.Range("N20").Value = -1181.3333
$ws.Range("H134").Value = 1804.2572
$ws.Range("I134").Value = 1924.24
$ws.Range("J134").Value = 1504.3
$ws.Range("K134").Value = 5772.72
$ws.Range("L134").Value = 4512.9
$ws.Range("M134").Value = -3237.72
$ws.Range("N134").Value = -9582.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3147.1785
$ws.Range("I134").Value = 2167.8655
$ws.Range("K134").Value = 6503.5965
$ws.Range("M134").Value = -3968.5965

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1344.2222
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 1387.25
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 4161.75
$ws.Range("M113").Value = -830
$ws.Range("N113").Value = -8501.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7048.273
$ws.Range("I70").Value = 6896.4287
$ws.Range("J70").Value = 7314
$ws.Range("K70").Value = 6896.4287
$ws.Range("L70").Value = 7314
$ws.Range("M70").Value = -6626.4287
$ws.Range("N70").Value = -7854
$ws.Range("H73").Value = 7048.273
$ws.Range("I73").Value = 6896.4287
$ws.Range("J73").Value = 7314
$ws.Range("K73").Value = 6896.4287
$ws.Range("L73").Value = 7314
$ws.Range("M73").Value = -5960.4287
$ws.Range("N73").Value = -9186
$ws.Range("H107").Value = 912.2174
$ws.Range("I107").Value = 351.16666
$ws.Range("J107").Value = 1524.2727
$ws.Range("K107").Value = 351.16666
$ws.Range("L107").Value = 1524.2727
$ws.Range("M107").Value = 1568.83334
$ws.Range("N107").Value = -5364.2727
$ws.Range("H122").Value = 1657.6364
$ws.Range("I122").Value = 1620.5714
$ws.Range("J122").Value = 1722.5
$ws.Range("K122").Value = 4861.7142
$ws.Range("L122").Value = 5167.5
$ws.Range("M122").Value = -2411.7142
$ws.Range("N122").Value = -10067.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3512.6924
$ws.Range("J16").Value = 9999.75
$ws.Range("L16").Value = 9999.75
$ws.Range("N16").Value = -10339.75
$ws.Range("H22").Value = 1522.375
$ws.Range("I22").Value = 1264.1428
$ws.Range("K22").Value = 1264.1428
$ws.Range("M22").Value = -969.1428000000001
$ws.Range("H27").Value = 1522.375
$ws.Range("I27").Value = 1264.1428
$ws.Range("K27").Value = 1264.1428
$ws.Range("M27").Value = -1157.1428
$ws.Range("H55").Value = 469
$ws.Range("I55").Value = 381.7
$ws.Range("J55").Value = 687.25
$ws.Range("K55").Value = 381.7
$ws.Range("L55").Value = 687.25
$ws.Range("M55").Value = -208.7
$ws.Range("N55").Value = -1033.25
$ws.Range("H100").Value = 2166.9614
$ws.Range("I100").Value = 1416
$ws.Range("J100").Value = 2810.6428
$ws.Range("K100").Value = 1416
$ws.Range("L100").Value = 2810.6428
$ws.Range("M100").Value = -875
$ws.Range("N100").Value = -3892.6428
$ws.Range("H132").Value = 1958.2333
$ws.Range("I132").Value = 1664.6666
$ws.Range("K132").Value = 4993.9998
$ws.Range("M132").Value = -2463.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1491.8182
$ws.Range("I81").Value = 1182.625
$ws.Range("J81").Value = 2316.3333
$ws.Range("K81").Value = 2365.25
$ws.Range("L81").Value = 4632.6666
$ws.Range("M81").Value = -1304.25
$ws.Range("N81").Value = -6754.6666
$ws.Range("H84").Value = 1491.8182
$ws.Range("I84").Value = 1182.625
$ws.Range("J84").Value = 2316.3333
$ws.Range("K84").Value = 11826.25
$ws.Range("L84").Value = 23163.333
$ws.Range("M84").Value = -6522.25
$ws.Range("N84").Value = -33771.333
$ws.Range("H113").Value = 703.06665
$ws.Range("I113").Value = 561.63635
$ws.Range("J113").Value = 1092
$ws.Range("K113").Value = 1684.90905
$ws.Range("L113").Value = 3276
$ws.Range("M113").Value = 485.09095
$ws.Range("N113").Value = -7616
$ws.Range("H122").Value = 3269.6597
$ws.Range("I122").Value = 3172.9111
$ws.Range("K122").Value = 9518.7333
$ws.Range("M122").Value = -7068.7333
$ws.Range("H130").Value = 49000
$ws.Range("J130").Value = 49000
$ws.Range("L130").Value = 49000
$ws.Range("N130").Value = -59040
$ws.Range("H132").Value = 3196
$ws.Range("I132").Value = 2476.5518
$ws.Range("J132").Value = 7368.8
$ws.Range("K132").Value = 7429.655400000001
$ws.Range("L132").Value = 22106.4
$ws.Range("M132").Value = -4899.655400000001
$ws.Range("N132").Value = -27166.4
